$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.429.38"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.866.46"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7049"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3142"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07854"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08027"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").Value = "1.887.86"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.194"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7009"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.450"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "29.501.85"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008328"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.36%  "
$ws.Range("D20").Value = "2.141.27"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.600"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1553"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.016"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.501"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.324"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.261"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05298"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.886"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7467"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").Value = "1.259.67"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.745"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8961"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.947"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000128"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "2.039.37"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5190"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.790"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.489"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4306"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.82%  "
